$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-10 Thursday", "2024-10-11 Friday"),
    @("81×43=", "73×69="),
    @("61×73=", "66×27="),
    @("47×89=", "13×28="),
    @("27×64=", "12×49="),
    @("54×58=", "45×84="),
    @("94×90=", "58×27="),
    @("46×33=", "20×16="),
    @("94×44=", "29×97="),
    @("72×51=", "76×20="),
    @("23×13=", "34×19="),
    @("41×89=", "19×14="),
    @("90×67=", "68×18="),
    @("39×46=", "95×76="),
    @("73×92=", "38×30="),
    @("59×95=", "25×44="),
    @("70×46=", "14×27="),
    @("65×51=", "70×30="),
    @("29×41=", "25×67="),
    @("15×46=", "70×31="),
    @("40×51=", "90×57="),
    @("83×85=", "20×98="),
    @("83×49=", "30×93="),
    @("69×11=", "99×89="),
    @("38×63=", "38×68="),
    @("77×97=", "97×14=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
